$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8363156171869264
$ws.Range("C2").Value = 0.1117734916592923
$ws.Range("D2").Value = 0.07783525045276463
$ws.Range("E2").Value = 0.1024704006889756
$ws.Range("G2").Value = 0.002491063255316643
$ws.Range("K2").Value = 0.4541948751577252
$ws.Range("L2").Value = 0.2043373006987821
$ws.Range("N2").Value = 2.087643721067938
$ws.Range("O2").Value = 4.667196951816209
$ws.Range("B3").Value = 0.7930457608114807
$ws.Range("C3").Value = 0.1106515677999624
$ws.Range("D3").Value = 0.07070596847079003
$ws.Range("E3").Value = 0.1018676378114058
$ws.Range("G3").Value = 0.002493974534854532
$ws.Range("K3").Value = 0.4155658996021145
$ws.Range("L3").Value = 0.197327879259106
$ws.Range("N3").Value = 2.107374676167275
$ws.Range("O3").Value = 4.673320793562823
$ws.Range("B4").Value = 0.7668591186546792
$ws.Range("C4").Value = 0.1099532293252281
$ws.Range("D4").Value = 0.06636338544122111
$ws.Range("E4").Value = 0.1015473777051668
$ws.Range("G4").Value = 0.00249585844233875
$ws.Range("K4").Value = 0.3919950610853959
$ws.Range("L4").Value = 0.1931363442196528
$ws.Range("N4").Value = 2.120133075942984
$ws.Range("O4").Value = 4.679344104076307
$ws.Range("B5").Value = 0.7562842168974555
$ws.Range("C5").Value = 0.1096662779686355
$ws.Range("D5").Value = 0.06460250007175716
$ws.Range("E5").Value = 0.1014294161733886
$ws.Range("G5").Value = 0.002496650456142199
$ws.Range("K5").Value = 0.3824272049627666
$ws.Range("L5").Value = 0.1914565384443279
$ws.Range("N5").Value = 2.125494014435709
$ws.Range("O5").Value = 4.682367818351679
$ws.Range("B6").Value = 0.754534097803969
$ws.Range("C6").Value = 0.1096184868787589
$ws.Range("D6").Value = 0.064310634898888
$ws.Range("E6").Value = 0.1014105869972219
$ws.Range("G6").Value = 0.002496783439427212
$ws.Range("K6").Value = 0.3808407423757956
$ws.Range("L6").Value = 0.1911793172945409
$ws.Range("N6").Value = 2.126393964899123
$ws.Range("O6").Value = 4.682904281780651
$ws.Range("B7").Value = 0.7667161109718847
$ws.Range("C7").Value = 0.1099493689946627
$ws.Range("D7").Value = 0.06633960209191514
$ws.Range("E7").Value = 0.1015457360123797
$ws.Range("G7").Value = 0.002495869025285553
$ws.Range("K7").Value = 0.3918658735122165
$ws.Range("L7").Value = 0.193113575210603
$ws.Range("N7").Value = 2.120204720407674
$ws.Range("O7").Value = 4.679382578347031
$ws.Range("B8").Value = 0.8213173354588434
$ws.Range("C8").Value = 0.1113886235954382
$ws.Range("D8").Value = 0.07536984118594603
$ws.Range("E8").Value = 0.1022522323942461
$ws.Range("G8").Value = 0.00249204710703177
$ws.Range("K8").Value = 0.4408451921290464
$ws.Range("L8").Value = 0.2018971735254098
$ws.Range("N8").Value = 2.094313376982566
$ws.Range("O8").Value = 4.668838737293044
$ws.Range("B9").Value = 0.9313990331465334
$ws.Range("C9").Value = 0.1141355475589734
$ws.Range("D9").Value = 0.09335558895557483
$ws.Range("E9").Value = 0.1040326659131665
$ws.Range("G9").Value = 0.002485313654592888
$ws.Range("K9").Value = 0.5380532003003964
$ws.Range("L9").Value = 0.220011791889192
$ws.Range("N9").Value = 2.048647125427024
$ws.Range("O9").Value = 4.666122741748751
$ws.Range("B10").Value = 1.014098052923998
$ws.Range("C10").Value = 0.1161075012547528
$ws.Range("D10").Value = 0.1067418457980409
$ws.Range("E10").Value = 0.1055812311562043
$ws.Range("G10").Value = 0.002480826033252774
$ws.Range("K10").Value = 0.6101720689757144
$ws.Range("L10").Value = 0.2338636877011027
$ws.Range("N10").Value = 2.018209023833151
$ws.Range("O10").Value = 4.675087370947722
$ws.Range("B11").Value = 1.052113575151054
$ws.Range("C11").Value = 0.1169945250511617
$ws.Range("D11").Value = 0.1128696878230784
$ws.Range("E11").Value = 0.1063378830838566
$ws.Range("G11").Value = 0.002478883268352783
$ws.Range("K11").Value = 0.6431318509854975
$ws.Range("L11").Value = 0.240283412380677
$ws.Range("N11").Value = 2.005037845377107
$ws.Range("O11").Value = 4.681548396607411
$ws.Range("B12").Value = 1.066565533051062
$ws.Range("C12").Value = 0.1173289689393116
$ws.Range("D12").Value = 0.1151956841728463
$ws.Range("E12").Value = 0.1066319040977213
$ws.Range("G12").Value = 0.002478161707961562
$ws.Range("K12").Value = 0.6556345420018488
$ws.Range("L12").Value = 0.2427314045077935
$ws.Range("N12").Value = 2.000147408409763
$ws.Range("O12").Value = 4.684337805304551
$ws.Range("B13").Value = 1.063450549806191
$ws.Range("C13").Value = 0.1172570052473745
$ws.Range("D13").Value = 0.1146944939655015
$ws.Range("E13").Value = 0.106568248353458
$ws.Range("G13").Value = 0.002478316481814041
$ws.Range("K13").Value = 0.6529409122864251
$ws.Range("L13").Value = 0.2422034311043859
$ws.Range("N13").Value = 2.0011963267807
$ws.Range("O13").Value = 4.683721810490908
$ws.Range("B14").Value = 1.053301421101082
$ws.Range("C14").Value = 0.1170220691291775
$ws.Range("D14").Value = 0.1130609384131844
$ws.Range("E14").Value = 0.1063619222674248
$ws.Range("G14").Value = 0.002478823622329494
$ws.Range("K14").Value = 0.6441600255554363
$ws.Range("L14").Value = 0.2404844701033539
$ws.Range("N14").Value = 2.004633557519082
$ws.Range("O14").Value = 4.681771013023848
$ws.Range("B15").Value = 1.047092103157581
$ws.Range("C15").Value = 0.1168779744587525
$ws.Range("D15").Value = 0.1120610567041496
$ws.Range("E15").Value = 0.1062365170453639
$ws.Range("G15").Value = 0.002479136097850057
$ws.Range("K15").Value = 0.6387842717812475
$ws.Range("L15").Value = 0.2394337669467603
$ws.Range("N15").Value = 2.006751622243254
$ws.Range("O15").Value = 4.680620733239948
$ws.Range("B16").Value = 1.0116215576939
$ws.Range("C16").Value = 0.116049329455123
$ws.Range("D16").Value = 0.106342148830862
$ws.Range("E16").Value = 0.105532831478655
$ws.Range("G16").Value = 0.002480954977488957
$ws.Range("K16").Value = 0.6080211052189384
$ws.Range("L16").Value = 0.2334465225634261
$ws.Range("N16").Value = 2.019083387883743
$ws.Range("O16").Value = 4.674713086844179
$ws.Range("B17").Value = 0.9899623928012886
$ws.Range("C17").Value = 0.1155384065266034
$ws.Range("D17").Value = 0.1028436098175831
$ws.Range("E17").Value = 0.105114503922568
$ws.Range("G17").Value = 0.002482096026869032
$ws.Range("K17").Value = 0.589187685863692
$ws.Range("L17").Value = 0.2298038401443421
$ws.Range("N17").Value = 2.026821548786653
$ws.Range("O17").Value = 4.671699336532612
$ws.Range("B18").Value = 0.9775418561316371
$ws.Range("C18").Value = 0.1152435939455358
$ws.Range("D18").Value = 0.1008349540907147
$ws.Range("E18").Value = 0.1048788072514562
$ws.Range("G18").Value = 0.002482761620056131
$ws.Range("K18").Value = 0.5783695834775244
$ws.Range("L18").Value = 0.2277198154385331
$ws.Range("N18").Value = 2.031335887370698
$ws.Range("O18").Value = 4.670190197973454
$ws.Range("B19").Value = 0.9733428887270748
$ws.Range("C19").Value = 0.1151436138013153
$ws.Range("D19").Value = 0.1001554780689133
$ws.Range("E19").Value = 0.1047998488420276
$ws.Range("G19").Value = 0.002482988576586844
$ws.Range("K19").Value = 0.5747092422265325
$ws.Range("L19").Value = 0.2270161169425791
$ws.Range("N19").Value = 2.032875277166919
$ws.Range("O19").Value = 4.669717749075403
$ws.Range("B20").Value = 0.9922641973900852
$ws.Range("C20").Value = 0.1155928928419883
$ws.Range("D20").Value = 0.1032156615130333
$ws.Range("E20").Value = 0.1051585270918629
$ws.Range("G20").Value = 0.002481973598981199
$ws.Range("K20").Value = 0.5911910498832924
$ws.Range("L20").Value = 0.2301904562842765
$ws.Range("N20").Value = 2.02599123012175
$ws.Range("O20").Value = 4.671996941474816
$ws.Range("B21").Value = 1.056280942449519
$ws.Range("C21").Value = 0.1170911150754037
$ws.Range("D21").Value = 0.1135406036008533
$ws.Range("E21").Value = 0.1064223219782683
$ws.Range("G21").Value = 0.002478674280143496
$ws.Range("K21").Value = 0.6467386035162122
$ws.Range("L21").Value = 0.2409889101232068
$ws.Range("N21").Value = 2.003621320672639
$ws.Range("O21").Value = 4.682334706428605
$ws.Range("B22").Value = 1.098447521205514
$ws.Range("C22").Value = 0.1180618170759189
$ws.Range("D22").Value = 0.1203206979393769
$ws.Range("E22").Value = 0.10729195583475
$ws.Range("G22").Value = 0.002476600272184712
$ws.Range("K22").Value = 0.6831674891352861
$ws.Range("L22").Value = 0.2481452839168981
$ws.Range("N22").Value = 1.989567910617808
$ws.Range("O22").Value = 4.691088935947334
$ws.Range("B23").Value = 1.075912597610625
$ws.Range("C23").Value = 0.1175445139628053
$ws.Range("D23").Value = 0.1166990931145904
$ws.Range("E23").Value = 0.1068238243618751
$ws.Range("G23").Value = 0.002477699701988406
$ws.Range("K23").Value = 0.663713369810921
$ws.Range("L23").Value = 0.2443167536353599
$ws.Range("N23").Value = 1.997016600555288
$ws.Range("O23").Value = 4.686233800246356
$ws.Range("B24").Value = 0.9912234530469561
$ws.Range("C24").Value = 0.1155682629304025
$ws.Range("D24").Value = 0.1030474485301198
$ws.Range("E24").Value = 0.1051386092480939
$ws.Range("G24").Value = 0.002482028918872178
$ws.Range("K24").Value = 0.5902852994482259
$ws.Range("L24").Value = 0.2300156353477405
$ws.Range("N24").Value = 2.026366412952591
$ws.Range("O24").Value = 4.671861698063282
$ws.Range("B25").Value = 0.9012981695879034
$ws.Range("C25").Value = 0.1134005317208491
$ws.Range("D25").Value = 0.08845995276469409
$ws.Range("E25").Value = 0.1035087574487115
$ws.Range("G25").Value = 0.002487054212762218
$ws.Range("K25").Value = 0.5116324117799991
$ws.Range("L25").Value = 0.2150159724582039
$ws.Range("N25").Value = 2.060454395862454
$ws.Range("O25").Value = 4.683921810490908
